$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-15 Friday", "2024-11-16 Saturday"),
    @("244÷8=", "290÷9="),
    @("681÷3=", "246÷8="),
    @("119÷7=", "713÷2="),
    @("914÷9=", "976÷3="),
    @("344÷4=", "217÷6="),
    @("924÷2=", "611÷5="),
    @("103÷4=", "682÷2="),
    @("318÷9=", "559÷4="),
    @("477÷3=", "766÷9="),
    @("373÷6=", "837÷6="),
    @("124÷4=", "274÷9="),
    @("336÷6=", "745÷8="),
    @("255÷6=", "225÷6="),
    @("271÷2=", "838÷4="),
    @("774÷2=", "178÷6="),
    @("623÷9=", "491÷3="),
    @("479÷7=", "816÷9="),
    @("296÷7=", "449÷2="),
    @("922÷9=", "651÷3="),
    @("915÷8=", "897÷7="),
    @("649÷2=", "496÷7="),
    @("613÷5=", "221÷6="),
    @("995÷3=", "565÷9="),
    @("743÷2=", "889÷7="),
    @("610÷5=", "998÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Replacements complete"
